$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: make D2:D7 a shared formula (same formula as before, now entered as one range) ---
$ws.Range("D2:D7").Formula = '=C2*$H$2/B2'

# --- New columns K (x) and L (y): probability -> reward-share polynomial ---
$ws.Range("K1").Value = "x"
$ws.Range("L1").Value = "y"

$ws.Range("K2").Value = 0
$ws.Range("L2").Formula = '=924*POWER(K2,13) - 6006 * POWER(K2,12) + 16380 * POWER(K2,11) - 24024 * POWER(K2,10)  + 20020 * POWER(K2,9) - 9009 * POWER(K2,8) + 1716 * POWER(K2,7)'

$kValues = @(0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = 3 + $i
    $ws.Range("K$row").Value = $kValues[$i]
}
$ws.Range("L3:L12").Formula = '=924*POWER(K3,13) - 6006 * POWER(K3,12) + 16380 * POWER(K3,11) - 24024 * POWER(K3,10)  + 20020 * POWER(K3,9) - 9009 * POWER(K3,8) + 1716 * POWER(K3,7)'

# --- sheet view bookkeeping ---
$ws.Range("O15").Select() | Out-Null

$excel.CalculateFullRebuild()
